# Apply the Vupiter Gantt chart update:
#  - Bump progress on several existing tasks
#  - Extend/retime several phase-2 tasks
#  - Insert a new "Project Proposal 2" task row after "Project Proposal 1" (formerly "Diposition")
#  - Rename the renamed phase/tasks rows ("Diposition"/"ALL" -> "Project Proposal 1"/"Group", etc.)
#  - Replace the placeholder "Build Device" phase + Task 1-5 rows with the real
#    "Capstone 2 (Build)" phase and its tasks
#  - Move the active selection to D16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# ---- Phase "Trade Studies" tasks: bump progress to 100% ----
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 1
$ws.Range("D9").Value = 1
$ws.Range("D10").Value = 1
$ws.Range("D11").Value = 1

# Row 10 task renamed Display -> User Interface
$ws.Range("B10").Value = "User Interface"

# Extend end dates / durations for rows 9 and 11 (row 10 shares row 9's formula)
$ws.Range("F9").Formula = "=E9+27"
$ws.Range("F11").Formula = "=E11+27"

# ---- Phase "Design Specs" tasks: update progress + extend end dates ----
$ws.Range("D13").Value = 0.5
$ws.Range("D14").Value = 0.25
$ws.Range("D15").Value = 0.25
$ws.Range("D16").Value = 0.75

$ws.Range("F13").Formula = "=E13+20"
$ws.Range("F14").Formula = "=E14+20"

# Row 17/18 (Technical Approach / Project Description) retime
$ws.Range("E17").Formula = "=E16+7"
$ws.Range("F17").Formula = "=E17+13"
$ws.Range("F18").Formula = "=E18+13"

# ---- Phase "Management Approch": Timeline progress ----
$ws.Range("D21").Value = 0.75

# ---- Rename "Diposition"/"ALL" task to "Project Proposal 1"/"Group" and retime ----
$ws.Range("B23").Value = "Project Proposal 1"
$ws.Range("C23").Value = "Group"
$ws.Range("E23").Formula = "=E21"
$ws.Range("F23").Formula = "=E23+13"
$ws.Range("H23").ClearContents()

# ---- Insert a new task row 24: "Project Proposal 2" / "Group" ----
$ws.Range("A23:BL23").Copy()
$ws.Rows("24:24").Insert()
$ws.Range("B24").Value = "Project Proposal 2"
$ws.Range("C24").Value = "Group"
$ws.Range("D24").Value = 0
$ws.Range("E24").Formula = "=F23+1"
$ws.Range("F24").Formula = "=E24+17"

# ---- Rename the old "Build Device" phase header (now row 25) ----
$ws.Range("B25").Value = "Capstone 2 (Build)"

# ---- Replace the 5 placeholder Task N rows (now rows 26-30) with real tasks ----
$ws.Range("B26").Value = "Identify source, order parts"
$ws.Range("C26").Value = "Group"
$ws.Range("E26").Value = 44207
$ws.Range("F26").Formula = "=E26+13"

$ws.Range("B27").Value = "Build Vupiter"
$ws.Range("C27").Value = "Group"
$ws.Range("E27").Formula = "=F26+8"
$ws.Range("F27").Formula = "=E27+34"

$ws.Range("B28").Value = "Test/Troubleshoot"
$ws.Range("C28").Value = "Group"
$ws.Range("E28").Formula = "=F27+8"
$ws.Range("F28").Formula = "=E28+13"

$ws.Range("B29").Value = "Fix Bugs and Finalize"
$ws.Range("C29").Value = "Group"
$ws.Range("E29").Formula = "=F28+1"
$ws.Range("F29").Formula = "=E29+20"

$ws.Range("B30").Value = "Project Submission"
$ws.Range("C30").Value = "Group"
$ws.Range("E30").Formula = "=F29+1"
$ws.Range("F30").Formula = "=E30"

# ---- Move active selection to D16, as in the saved file ----
$ws.Range("D16").Select()
